# repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the freshly re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = -4
    3  = 5
    4  = -6
    5  = -2
    6  = -3
    8  = 1
    9  = 2
    10 = 3
    11 = 8
    12 = 1
    13 = -3
    14 = 2
    15 = 8
    16 = -4
    17 = 2
    20 = -2
    21 = -1
    22 = 4
    23 = -1
    24 = 6
    25 = -2
    27 = -1
    28 = -1
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
